$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 407
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 18
$ws.Range("H2").Value = 16
$ws.Range("I2").Value = 16
$ws.Range("K2").Value = 3007
$ws.Range("L2").Value = 286
$ws.Range("M2").Value = 2721
$ws.Range("N2").Value = 2721
$ws.Range("P2").Value = 45
$ws.Range("Q2").Value = -50
$ws.Range("R2").Value = 155
$ws.Range("S2").Value = -14
$ws.Range("T2").Value = 6
$ws.Range("U2").Value = -56
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 0.75
$ws.Range("X2").Value = 3.86
$ws.Range("Y2").Value = 0.58
$ws.Range("Z2").Value = 0.52
$ws.Range("AA2").Value = 10.51
$ws.Range("AB2").Value = 5925.77
$ws.Range("AC2").Value = 1745
$ws.Range("AD2").Value = 65.61
$ws.Range("AE2").Value = 302358
$ws.Range("AF2").Value = 0.38
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 1.31
$ws.Range("AI2").Value = 85.95
$ws.Range("AJ2").Value = 900000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 1791
$ws.Range("E3").Value = 60
$ws.Range("F3").Value = 60
$ws.Range("G3").Value = 148
$ws.Range("H3").Value = 114
$ws.Range("I3").Value = 114
$ws.Range("K3").Value = 3177
$ws.Range("L3").Value = 339
$ws.Range("M3").Value = 2838
$ws.Range("N3").Value = 2838
$ws.Range("P3").Value = 45
$ws.Range("Q3").Value = 197
$ws.Range("R3").Value = -231
$ws.Range("S3").Value = -12
$ws.Range("T3").Value = 10
$ws.Range("U3").Value = 187
$ws.Range("V3").Value = 3
$ws.Range("W3").Value = 3.36
$ws.Range("X3").Value = 6.34
$ws.Range("Y3").Value = 4.08
$ws.Range("Z3").Value = 3.67
$ws.Range("AA3").Value = 11.93
$ws.Range("AB3").Value = 6184.43
$ws.Range("AC3").Value = 12613
$ws.Range("AD3").Value = 11.61
$ws.Range("AE3").Value = 315367
$ws.Range("AF3").Value = 0.46
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 1.02
$ws.Range("AI3").Value = 11.89
$ws.Range("AJ3").Value = 900000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1795
$ws.Range("E4").Value = 32
$ws.Range("F4").Value = 32
$ws.Range("G4").Value = 228
$ws.Range("H4").Value = 173
$ws.Range("I4").Value = 173
$ws.Range("K4").Value = 3345
$ws.Range("L4").Value = 342
$ws.Range("M4").Value = 3003
$ws.Range("N4").Value = 3003
$ws.Range("P4").Value = 45
$ws.Range("Q4").Value = 23
$ws.Range("R4").Value = 107
$ws.Range("S4").Value = -13
$ws.Range("T4").Value = 1
$ws.Range("U4").Value = 22
$ws.Range("V4").Value = 3
$ws.Range("W4").Value = 1.76
$ws.Range("X4").Value = 9.66
$ws.Range("Y4").Value = 5.93
$ws.Range("Z4").Value = 5.31
$ws.Range("AA4").Value = 11.4
$ws.Range("AB4").Value = 6551.39
$ws.Range("AC4").Value = 19255
$ws.Range("AD4").Value = 8.91
$ws.Range("AE4").Value = 333631
$ws.Range("AF4").Value = 0.51
$ws.Range("AG4").Value = 2500
$ws.Range("AH4").Value = 1.46
$ws.Range("AI4").Value = 12.98
$ws.Range("AJ4").Value = 900000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 1824
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 60
$ws.Range("G5").Value = 114
$ws.Range("H5").Value = 91
$ws.Range("I5").Value = 91
$ws.Range("K5").Value = 3316
$ws.Range("L5").Value = 268
$ws.Range("M5").Value = 3048
$ws.Range("N5").Value = 3048
$ws.Range("P5").Value = 45
$ws.Range("Q5").Value = 9
$ws.Range("R5").Value = 78
$ws.Range("S5").Value = -21
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = 5
$ws.Range("V5").Value = 4
$ws.Range("W5").Value = 3.27
$ws.Range("X5").Value = 5
$ws.Range("Y5").Value = 3.02
$ws.Range("Z5").Value = 2.74
$ws.Range("AA5").Value = 8.81
$ws.Range("AB5").Value = 6677.96
$ws.Range("AC5").Value = 10140
$ws.Range("AD5").Value = 14.79
$ws.Range("AE5").Value = 338648
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 1500
$ws.Range("AH5").Value = 1
$ws.Range("AI5").Value = 14.79
$ws.Range("AJ5").Value = 900000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 1714
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 27
$ws.Range("I6").Value = 27
$ws.Range("K6").Value = 3357
$ws.Range("L6").Value = 274
$ws.Range("M6").Value = 3083
$ws.Range("N6").Value = 3083
$ws.Range("P6").Value = 45
$ws.Range("Q6").Value = 45
$ws.Range("R6").Value = -1
$ws.Range("S6").Value = -16
$ws.Range("T6").Value = 18
$ws.Range("U6").Value = 26
$ws.Range("V6").Value = 2
$ws.Range("W6").Value = 0.44
$ws.Range("X6").Value = 1.58
$ws.Range("Y6").Value = 0.88
$ws.Range("Z6").Value = 0.81
$ws.Range("AA6").Value = 8.89
$ws.Range("AB6").Value = 6700.06
$ws.Range("AC6").Value = 3010
$ws.Range("AD6").Value = 43.52
$ws.Range("AE6").Value = 342569
$ws.Range("AF6").Value = 0.38
$ws.Range("AG6").Value = 1500
$ws.Range("AH6").Value = 1.15
$ws.Range("AI6").Value = 49.83
$ws.Range("AJ6").Value = 900000

# Row 7 - clear all data columns, keep A/B/C
$ws.Range("D7:AJ7").ClearContents()

# Row 8 - clear all data columns, keep A/B/C
$ws.Range("D8:AJ8").ClearContents()

# Row 9 - clear all data columns, keep A/B/C
$ws.Range("D9:AJ9").ClearContents()
